$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values to update: maps cell address -> new text value.
# Values that look like plain numbers (e.g. "1.00", "0.190") must be
# forced to Text so Excel does not coerce them (losing the trailing
# zeros / exact text form) the way the source data represents them.
$updates = [ordered]@{
    'D2' = '98.504.57'
    'E2' = '  -0.33%  '
    'D3' = '3.371.02'
    'E3' = '  -0.10%  '
    'E4' = '  -0.02%  '
    'D5' = '258.15'
    'E5' = '  -0.49%  '
    'D6' = '670.17'
    'E6' = '  +6.50%  '
    'D7' = '1.55'
    'E7' = '  +11.67%  '
    'E8' = '  +17.09%  '
    'D9' = '1.09'
    'E9' = '  +26.43%  '
    'E10' = '  +0.00%  '
    'D11' = '3.370.68'
    'E11' = '  -0.07%  '
    'E12' = '  +5.70%  '
    'D13' = '42.47'
    'E13' = '  +17.08%  '
    'D14' = '0.0000269'
    'E14' = '  +8.19%  '
    'D15' = '98.327.27'
    'E15' = '  -0.27%  '
    'B16' = 'WrappedliquidstakedEther2.0'
    'C16' = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
    'D16' = '4.007.77'
    'E16' = '  +0.17%  '
    'B17' = 'Toncoin'
    'C17' = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
    'D17' = '5.63'
    'E17' = '  +2.54%  '
    'B18' = 'WrappedEther'
    'C18' = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
    'D18' = '3.374.60'
    'E18' = '  -0.10%  '
    'B19' = 'Polkadot'
    'C19' = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
    'D19' = '7.67'
    'E19' = '  +25.14%  '
    'B20' = 'Chainlink'
    'C20' = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
    'D20' = '16.92'
    'E20' = '  +10.96%  '
    'E21' = '  +0.64%  '
    'B22' = 'BitcoinCash'
    'C22' = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
    'D22' = '530.59'
    'E22' = '  +8.18%  '
    'B23' = 'Uniswap'
    'C23' = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
    'D23' = '10.61'
    'E23' = '  +13.22%  '
    'B24' = 'Stellar'
    'C24' = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    'D24' = '0.449'
    'E24' = '  +57.64%  '
    'B25' = 'PEPE'
    'C25' = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
    'D25' = '0.0000213'
    'E25' = '  +0.75%  '
    'B26' = 'Litecoin'
    'C26' = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
    'D26' = '103.25'
    'E26' = '  +16.27%  '
    'B27' = 'NEARProtocol'
    'C27' = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    'D27' = '6.29'
    'E27' = '  +11.41%  '
    'B28' = 'Aptos'
    'C28' = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
    'D28' = '12.71'
    'E28' = '  +6.25%  '
    'B29' = 'WrappedeETH'
    'C29' = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
    'D29' = '3.550.93'
    'E29' = '  -0.07%  '
    'B30' = 'Hedera'
    'C30' = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
    'D30' = '0.150'
    'E30' = '  +10.68%  '
    'B31' = 'Dai'
    'C31' = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
    'D31' = '0.998'
    'E31' = '  -0.29%  '
    'B32' = 'InternetComputer(DFINITY)'
    'C32' = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
    'D32' = '11.18'
    'E32' = '  +15.24%  '
    'B33' = 'Cronos'
    'C33' = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
    'D33' = '0.190'
    'E33' = '  -0.79%  '
    'B34' = 'Binance-PegBSC-USD'
    'C34' = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
    'D34' = '1.00'
    'E34' = '  +0.39%  '
    'B35' = 'EthereumClassic'
    'C35' = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
    'D35' = '29.82'
    'E35' = '  +6.39%  '
    'B36' = 'PolygonEcosystemToken'
    'C36' = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
    'D36' = '0.541'
    'E36' = '  +17.43%  '
    'B37' = 'RenderToken'
    'C37' = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
    'D37' = '7.91'
    'E37' = '  +8.09%  '
    'B38' = 'Kaspa'
    'C38' = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
    'D38' = '0.162'
    'E38' = '  +7.08%  '
    'B39' = 'PancakeSwap'
    'C39' = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
    'D39' = '2.14'
    'E39' = '  +8.90%  '
    'B40' = 'Bittensor'
    'C40' = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
    'D40' = '529.05'
    'E40' = '  +5.67%  '
    'B41' = 'VeChain'
    'C41' = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    'D41' = '0.0453'
    'E41' = '  +38.57%  '
    'B42' = 'Fetch.AI'
    'C42' = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
    'D42' = '1.35'
    'E42' = '  +6.29%  '
    'B43' = 'WhiteBITCoin'
    'C43' = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
    'D43' = '24.71'
    'E43' = '  -0.86%  '
    'B44' = 'MantraDAO'
    'C44' = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
    'D44' = '3.78'
    'E44' = '  +0.71%  '
    'B45' = 'ARBITRUM'
    'C45' = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
    'D45' = '0.842'
    'E45' = '  +6.90%  '
    'B46' = 'dogwifhat'
    'C46' = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
    'D46' = '3.38'
    'E46' = '  +2.18%  '
    'B47' = 'USDe'
    'C47' = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
    'D47' = '1.00'
    'E47' = '  +0.01%  '
    'B48' = 'Cosmos'
    'C48' = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
    'D48' = '7.96'
    'E48' = '  +19.76%  '
    'D49' = '5.19'
    'E49' = '  +11.99%  '
    'B50' = 'Stacks'
    'C50' = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
    'D50' = '2.07'
    'E50' = '  +6.51%  '
    'D51' = '1.54'
    'E51' = '  +12.64%  '
}

foreach ($addr in $updates.Keys) {
    $value = $updates[$addr]
    $cell = $ws.Range($addr)
    if ($value -match '^[+-]?\d+(\.\d+)?$') {
        # Force text storage, then strip the formatting change back off
        # so the cell keeps its original (unstyled) appearance.
        $cell.NumberFormat = "@"
        $cell.Value = $value
        $cell.ClearFormats()
    } else {
        $cell.Value = $value
    }
}
